$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 8.382531
$ws.Range("H2").Value = 25.147593
$ws.Range("I2").Value = 0.3278601051951505
$ws.Range("J2").Value = 0.3278601051951506
$ws.Range("M2").Value = 0.6746743333333334
$ws.Range("N2").Value = 2.024023
$ws.Range("O2").Value = 0.07069047851636343
$ws.Range("P2").Value = 0.07069047851636343
$ws.Range("Q2").Value = 5.655478514071
$ws.Range("R2").Value = 50.899306626639
$ws.Range("S2").Value = 0.02317658772267044
$ws.Range("T2").Value = 0.02317658772267044

# Row 3
$ws.Range("G3").Value = 8.382531
$ws.Range("H3").Value = 25.147593
$ws.Range("I3").Value = 0.3278601051951505
$ws.Range("J3").Value = 0.3278601051951506
$ws.Range("O3").Value = 0.4692497642600617
$ws.Range("P3").Value = 0.4692497642600616
$ws.Range("Q3").Value = 37.541575827519
$ws.Range("R3").Value = 337.874182447671
$ws.Range("S3").Value = 0.1538482770731034
$ws.Range("T3").Value = 0.1538482770731034

# Row 4
$ws.Range("G4").Value = 8.382531
$ws.Range("H4").Value = 25.147593
$ws.Range("I4").Value = 0.3278601051951505
$ws.Range("J4").Value = 0.3278601051951506
$ws.Range("O4").Value = 0.460059757223575
$ws.Range("P4").Value = 0.460059757223575
$ws.Range("Q4").Value = 36.806344033509
$ws.Range("R4").Value = 331.257096301581
$ws.Range("S4").Value = 0.1508352403993767
$ws.Range("T4").Value = 0.1508352403993767

# Row 5
$ws.Range("I5").Value = 0.2503004183517279
$ws.Range("J5").Value = 0.250300418351728
$ws.Range("M5").Value = 0.6746743333333334
$ws.Range("N5").Value = 2.024023
$ws.Range("O5").Value = 0.07069047851636343
$ws.Range("P5").Value = 0.07069047851636343
$ws.Range("Q5").Value = 4.317599535962445
$ws.Range("R5").Value = 38.858395823662
$ws.Range("S5").Value = 0.0176938563461296
$ws.Range("T5").Value = 0.01769385634612961

# Row 6
$ws.Range("I6").Value = 0.2503004183517279
$ws.Range("J6").Value = 0.250300418351728
$ws.Range("O6").Value = 0.4692497642600617
$ws.Range("P6").Value = 0.4692497642600616
$ws.Range("S6").Value = 0.1174534123057431
$ws.Range("T6").Value = 0.1174534123057432

# Row 7
$ws.Range("I7").Value = 0.2503004183517279
$ws.Range("J7").Value = 0.250300418351728
$ws.Range("O7").Value = 0.460059757223575
$ws.Range("P7").Value = 0.460059757223575
$ws.Range("S7").Value = 0.1151531496998552
$ws.Range("T7").Value = 0.1151531496998552

# Row 8
$ws.Range("I8").Value = 0.4218394764531215
$ws.Range("J8").Value = 0.4218394764531215
$ws.Range("M8").Value = 0.6746743333333334
$ws.Range("N8").Value = 2.024023
$ws.Range("O8").Value = 0.07069047851636343
$ws.Range("P8").Value = 0.07069047851636343
$ws.Range("Q8").Value = 7.276591624490445
$ws.Range("R8").Value = 65.489324620414
$ws.Range("S8").Value = 0.02982003444756338
$ws.Range("T8").Value = 0.02982003444756338

# Row 9
$ws.Range("I9").Value = 0.4218394764531215
$ws.Range("J9").Value = 0.4218394764531215
$ws.Range("O9").Value = 0.4692497642600617
$ws.Range("P9").Value = 0.4692497642600616
$ws.Range("S9").Value = 0.1979480748812151
$ws.Range("T9").Value = 0.1979480748812151

# Row 10
$ws.Range("I10").Value = 0.4218394764531215
$ws.Range("J10").Value = 0.4218394764531215
$ws.Range("O10").Value = 0.460059757223575
$ws.Range("P10").Value = 0.460059757223575
$ws.Range("S10").Value = 0.1940713671243431
$ws.Range("T10").Value = 0.1940713671243431
